$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "Record"
$ws.Range("B5").Value = "Balanço Geral"
$ws.Range("C5").Value = "Transporte"
$ws.Range("D5").Value = "2025-04-02T12:47"
$ws.Range("E5").Value = "Neutro"
$ws.Range("F5").Value = "Vans do setor C voltam a circular após reunião com representantes da prefeitura. Repórter *ao vivo*. Ontem, motoristas de vans do setor C pararam em frente à prefeitura e pediram para conversar direto com o prefeito Wladimir Garotinho. Alegam que falta uma parte do repasse. Prefeitura ficou de fazer correção e os permissionários vão enviar outro relatório. Pagamento semana que vem da primeira quinzena de março."

$ws.Range("A6").Value = "Record"
$ws.Range("B6").Value = "Balanço Geral"
$ws.Range("C6").Value = "Saúde"
$ws.Range("D6").Value = "2025-04-02T12:58"
$ws.Range("E6").Value = "Neutro"
$ws.Range("F6").Value = "Começou hoje a vacinação contra a influenza em cidades do Norte Fluminense. Repórter *ao vivo*. Imagens da Secretaria de Saúde de Campos. Informações de Macaé e São João da Barra."
